$d = $word.ActiveDocument

# Replacement 1
$d.Content.Find.Execute(
    "Uzakufumana umyalezo okukhumbuzayo ukuba ugqibe isifundo sakho. Ukuba awuwubonanga umyalezo kulungile! Ungabuyela kwi ParentText nanini na ukuzikhumbuza ngesifundo sakho. And if you miss it, it is also okay! You can always return to ParentText anytime to catch up on your lesson.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Uzakufumana umyalezo okukhumbuzayo ukuba ugqibe isifundo sakho. Ukuba awuwubonanga umyalezo kulungile! Ungabuyela kwi ParentText nanini na ukuzikhumbuza ngesifundo sakho. Kwaye ukuba ikuphosile, kulungile! Usenokubuyela kwi ParentText nangaliphi na ixesha ukuze uqhubele phambili nesifundo sakho.",
    2
) | Out-Null

# Replacement 2
$d.Content.Find.Execute(
    "Ulwazi lwakho lukhuselekile apha: Alukho ulwazi okuzokwabelwana ngalo ngaphandle kwemvume yakho okanye luthengiswe ukwenza inzuzo. Umyalezo owuthumelayo uvaliwe kwaye utshixelwe kwiseva ekhuselekileyo. The messages you send are encrypted and locked in a secure server. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ulwazi lwakho lukhuselekile apha: Alukho ulwazi okuzokwabelwana ngalo ngaphandle kwemvume yakho okanye luthengiswe ukwenza inzuzo. Umyalezo owuthumelayo uvaliwe kwaye utshixelwe kwiseva ekhuselekileyo. Le miyalezo oyithumelayo inoguqulelo oluntsokothileyo kwaye itshixelwe kwiseva ekhuselekileyo. ",
    2
) | Out-Null

# Replacement 3
$d.Content.Find.Execute(
    "Khumbula, xa kukho umntu okwazi ukufikelela kwifoni yakho xa ingatshixwanga angakwazi ukubona imiyalezo yakho. Xa ngamanye amaxesha uthumela ulwazi olunobuzaza kwaye loonto ikukhathaza cima yonke imilayezo kwifoni yakho. So, if you send sensitive information and are worried, delete the messages from your phone. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Khumbula, xa kukho umntu okwazi ukufikelela kwifoni yakho xa ingatshixwanga angakwazi ukubona imiyalezo yakho. Xa ngamanye amaxesha uthumela ulwazi olunobuzaza kwaye loonto ikukhathaza cima yonke imilayezo kwifoni yakho. Ngoko ke, ukuba uthumela ulwazi ulunobuzaza kwaye unexhala, cima imiyalezo kwifowuni yakho. ",
    2
) | Out-Null

# Replacement 4
$d.Content.Find.Execute(
    "I-ParentText iya kubonelela ngamacebiso ngezifundo ezizakunceda ngobudlelwane bakho nomntwana wakho. Kukuwe ukusebenzisa lamacebiso uzame uwaprakthize. It is up to you to put these tips into practice!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I-ParentText iya kubonelela ngamacebiso ngezifundo ezizakunceda ngobudlelwane bakho nomntwana wakho. Kukuwe ukusebenzisa lamacebiso uzame uwaprakthize. Kuxhomekeke kuwe ukuba uwasebenzise la macebiso!",
    2
) | Out-Null

# Replacement 5
$d.Content.Find.Execute(
    "Enkosi kakhulu ukumamela! Ungayifumana le-vidiyo nanini na xa usiya kwi-Menu. Siyathemba ukonwabele ukuba kwi-ParentText nokuthi uzolusebenzisa ulwazi olufumene apha! You can access this video at any time via MENU. We hope you enjoy your ParentText journey and make the most out of it! ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enkosi kakhulu ukumamela! Ungayifumana le-vidiyo nanini na xa usiya kwi-Menu. Siyathemba ukonwabele ukuba kwi-ParentText nokuthi uzolusebenzisa ulwazi olufumene apha! Ungafikelela kulevidiyo nangaliphi na ixesha ngeMENYU. Siyathemba uya kukonwabela ukusebenzisa i ParentText kwaye wenze lukhulu kuyo! ",
    2
) | Out-Null
